$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 824.4666592499998
$schedule.Range("F2").Value = 13.63205455109127

# --- Detailed sheet updates ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B18").Value = 3.55269
$detailed.Range("B19").Value = 0.7
$detailed.Range("B20").Value = 0.009379999999999999
$detailed.Range("B21").Value = 0.0005
$detailed.Range("C21").Value = "historical"
$detailed.Range("B22").Value = 0
$detailed.Range("C22").Value = "historical"
$detailed.Range("B23").Value = -2.91592
$detailed.Range("B24").Value = -6.40473
$detailed.Range("B25").Value = -6.14838
$detailed.Range("B26").Value = -6.33429
$detailed.Range("B27").Value = -7.11669
$detailed.Range("B28").Value = -8.96598
$detailed.Range("B29").Value = -8.18792
$detailed.Range("B30").Value = -8.118080000000001
$detailed.Range("B31").Value = -7.54132
$detailed.Range("B32").Value = -6.46636
$detailed.Range("B33").Value = -5.26642
$detailed.Range("B35").Value = -0.9120200000000001
$detailed.Range("B36").Value = 0.00001
$detailed.Range("B37").Value = 22.05184
$detailed.Range("B38").Value = 30.06361
$detailed.Range("B39").Value = 47.37894
$detailed.Range("B40").Value = 60.85009
$detailed.Range("B41").Value = 57.22391
$detailed.Range("B42").Value = 62.22274
$detailed.Range("B43").Value = 57.59874
$detailed.Range("B44").Value = 58.85308
$detailed.Range("B45").Value = 57.3
$detailed.Range("B46").Value = 57.06007
$detailed.Range("B47").Value = 58.6455
$detailed.Range("B48").Value = 61.33335
$detailed.Range("B49").Value = 57.3
